$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 170.055555
$ws.Range("H2").Value = 510.166665
$ws.Range("I2").Value = 0.5874625966152389
$ws.Range("J2").Value = 0.587462596615239
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 45.95651366666667
$ws.Range("N2").Value = 137.869541
$ws.Range("O2").Value = 0.6189188856627118
$ws.Range("P2").Value = 0.6189188856627118
$ws.Range("Q2").Value = 7815.160437450085
$ws.Range("R2").Value = 70336.44393705076
$ws.Range("S2").Value = 0.3635916956656268
$ws.Range("T2").Value = 0.3635916956656269

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 170.055555
$ws.Range("H3").Value = 510.166665
$ws.Range("I3").Value = 0.5874625966152389
$ws.Range("J3").Value = 0.587462596615239
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 6.849914666666667
$ws.Range("N3").Value = 20.549744
$ws.Range("O3").Value = 0.09225115688993263
$ws.Range("P3").Value = 0.09225115688993261
$ws.Range("Q3").Value = 1164.86604034264
$ws.Range("R3").Value = 10483.79436308376
$ws.Range("S3").Value = 0.05419410416731961
$ws.Range("T3").Value = 0.05419410416731961

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 170.055555
$ws.Range("H4").Value = 510.166665
$ws.Range("I4").Value = 0.5874625966152389
$ws.Range("J4").Value = 0.587462596615239
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 21.446458
$ws.Range("N4").Value = 64.33937399999999
$ws.Range("O4").Value = 0.2888299574473556
$ws.Range("P4").Value = 0.2888299574473556
$ws.Range("Q4").Value = 3647.089317974189
$ws.Range("R4").Value = 32823.80386176771
$ws.Range("S4").Value = 0.1696767967822925
$ws.Range("T4").Value = 0.1696767967822925

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 69.750951
$ws.Range("H5").Value = 209.252853
$ws.Range("I5").Value = 0.2409569907365996
$ws.Range("J5").Value = 0.2409569907365995
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 45.95651366666667
$ws.Range("N5").Value = 137.869541
$ws.Range("O5").Value = 0.6189188856627118
$ws.Range("P5").Value = 0.6189188856627118
$ws.Range("Q5").Value = 3205.510532894497
$ws.Range("R5").Value = 28849.59479605047
$ws.Range("S5").Value = 0.1491328321993366
$ws.Range("T5").Value = 0.1491328321993366

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 69.750951
$ws.Range("H6").Value = 209.252853
$ws.Range("I6").Value = 0.2409569907365996
$ws.Range("J6").Value = 0.2409569907365995
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 6.849914666666667
$ws.Range("N6").Value = 20.549744
$ws.Range("O6").Value = 0.09225115688993263
$ws.Range("P6").Value = 0.09225115688993261
$ws.Range("Q6").Value = 477.788062268848
$ws.Range("R6").Value = 4300.092560419632
$ws.Range("S6").Value = 0.02222856115616809
$ws.Range("T6").Value = 0.02222856115616808

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 69.750951
$ws.Range("H7").Value = 209.252853
$ws.Range("I7").Value = 0.2409569907365996
$ws.Range("J7").Value = 0.2409569907365995
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 21.446458
$ws.Range("N7").Value = 64.33937399999999
$ws.Range("O7").Value = 0.2888299574473556
$ws.Range("P7").Value = 0.2888299574473556
$ws.Range("Q7").Value = 1495.910841081558
$ws.Range("R7").Value = 13463.19756973402
$ws.Range("S7").Value = 0.06959559738109491
$ws.Range("T7").Value = 0.0695955973810949

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 49.66818733333334
$ws.Range("H8").Value = 149.004562
$ws.Range("I8").Value = 0.1715804126481615
$ws.Range("J8").Value = 0.1715804126481615
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 45.95651366666667
$ws.Range("N8").Value = 137.869541
$ws.Range("O8").Value = 0.6189188856627118
$ws.Range("P8").Value = 0.6189188856627118
$ws.Range("Q8").Value = 2282.576729982894
$ws.Range("R8").Value = 20543.19056984604
$ws.Range("S8").Value = 0.1061943577977484
$ws.Range("T8").Value = 0.1061943577977484

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 49.66818733333334
$ws.Range("H9").Value = 149.004562
$ws.Range("I9").Value = 0.1715804126481615
$ws.Range("J9").Value = 0.1715804126481615
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 6.849914666666667
$ws.Range("N9").Value = 20.549744
$ws.Range("O9").Value = 0.09225115688993263
$ws.Range("P9").Value = 0.09225115688993261
$ws.Range("Q9").Value = 340.2228448813476
$ws.Range("R9").Value = 3062.005603932128
$ws.Range("S9").Value = 0.01582849156644493
$ws.Range("T9").Value = 0.01582849156644493

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 49.66818733333334
$ws.Range("H10").Value = 149.004562
$ws.Range("I10").Value = 0.1715804126481615
$ws.Range("J10").Value = 0.1715804126481615
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 21.446458
$ws.Range("N10").Value = 64.33937399999999
$ws.Range("O10").Value = 0.2888299574473556
$ws.Range("P10").Value = 0.2888299574473556
$ws.Range("Q10").Value = 1065.206693580465
$ws.Range("R10").Value = 9586.860242224187
$ws.Range("S10").Value = 0.04955756328396821
$ws.Range("T10").Value = 0.0495575632839682
